# This workbook tracks daily Limón (lemon) prices reported by
# "Femacal de La Calera". The edit adds a new week of price records:
# 4 brand-new data rows are inserted right before the existing row 1037,
# pushing all the rows that used to start at 1037 down by four (to 1041+).
# The worksheet's dimension grows from A1:T1112 to A1:T1116 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 1037; this shifts old rows 1037:1112
# down to 1041:1116 and auto-extends the used range / dimension.
$ws.Rows("1037:1040").Insert()

# Data for the 4 new rows (columns A-T), derived from the target diff.
$newRows = @(
    @{ Row = 1037; D = 44783; L = "1a amarillo"; M = 297; N = 2700; O = 3000; P = 2824; R = "Cabildo";               S = 176 },
    @{ Row = 1038; D = 44783; L = "1a amarillo"; M = 182; N = 2800; O = 3000; P = 2907; R = "Provincia de Quillota"; S = 182 },
    @{ Row = 1039; D = 44783; L = "2a amarillo"; M = 310; N = 2200; O = 2500; P = 2331; R = "Cabildo";               S = 146 },
    @{ Row = 1040; D = 44783; L = "2a amarillo"; M = 189; N = 2300; O = 2500; P = 2394; R = "Provincia de Quillota"; S = 150 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = 3
    $ws.Cells.Item($row, 2).Value  = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"

    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(1036, 4).NumberFormat

    $ws.Cells.Item($row, 5).Value  = 5
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100102
    $ws.Cells.Item($row, 8).Value  = "Cítricos"
    $ws.Cells.Item($row, 9).Value  = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/malla 16 kilos"
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
